$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.677.65'
$ws.Range("E2").Value = '  -0.27%  '

$ws.Range("D3").Value = '2.427.04'
$ws.Range("E3").Value = '  -3.12%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '486.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("E8").Value = '  +17.26%  '

$ws.Range("D9").Value = '2.427.01'
$ws.Range("E9").Value = '  -3.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0994'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("E11").Value = '  -0.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.125'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.16%  '

$ws.Range("D14").Value = '2.845.07'
$ws.Range("E14").Value = '  -3.02%  '

$ws.Range("D15").Value = '56.966.97'
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("E16").Value = '  -3.86%  '

$ws.Range("E17").Value = '  -3.35%  '

$ws.Range("D18").Value = '2.425.07'
$ws.Range("E18").Value = '  -3.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("E21").Value = '  -4.06%  '

$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '57.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.407'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.31%  '

$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("E27").Value = '  -2.36%  '

$ws.Range("D28").Value = '2.520.05'
$ws.Range("E28").Value = '  -3.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.12%  '

$ws.Range("D30").Value = '0.0₃0783'
$ws.Range("E30").Value = '  -4.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.51%  '

$ws.Range("E34").Value = '  -1.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.35'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.48%  '

$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("E37").Value = '  -2.73%  '

$ws.Range("E38").Value = '  -4.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.39%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '269.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("E45").Value = '  -4.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0532'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.80%  '

$ws.Range("D51").Value = '1.864.76'
$ws.Range("E51").Value = '  -2.51%  '
